# The Education section entries are stored as rich-text content controls
# (w:sdt) that are still showing their gallery placeholder text (w:showingPlcHdr).
# Replace each placeholder with the real text and let Word drop the now
# unused content control, leaving a plain run behind - exactly like a user
# clicking into the placeholder and typing over it.

$d = $word.ActiveDocument

foreach ($cc in $d.ContentControls) {
    if ($cc.ID -eq -1770613245) {
        # "MBA | Jasper University" -> "MBA | Thammasat University"
        $r = $cc.Range
        $r.Text = "MBA | Thammasat University"
        $cc.Delete()
    }
    elseif ($cc.ID -eq -1309927247) {
        # "BS Computer Science | Bellows College" -> "BS Computer Science | Chulalongkorn Univeristy"
        $r = $cc.Range
        $r.Text = "BS Computer Science | Chulalongkorn Univeristy"
        $cc.Delete()
    }
}
